$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of marks for "Janifer"
$ws.Range("A8").Value = "Janifer"
$ws.Range("B8").Value = 23
$ws.Range("C8").Value = 45
$ws.Range("D8").Value = 53

# Update selection to match where Excel would leave the cursor after entry
$ws.Range("D9").Select()
